# Applies the pending genome-annotation updates for rows 132-140 that were
# missed in the previous commit (new Transporters/Nitrogen/Sulfur/Motility/
# Other annotations), plus restores the sheet scroll position/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('J132').Value = 'MurNAc, chitobiose, glucose, glycolate, starch/glycogen'
$ws.Range('N132').Value = 'ammonia_assimilation, nitrogen_fixation, nitrous oxide reductase, one nitrate reductase'
$ws.Range('P132').Value = 'thiosulfate, one sulfite reductase, no evidence for sulfate oxidation'
$ws.Range('Q132').Value = 'branched amino, amino acid/amide, iron, LPS export, molybdate, oligopeptide/dipeptide, phosphate, phospholipid/cholesterol, sodium, tungstate/molybdate'
$ws.Range('S132').Value = 'chemotaxis (purine?), flagellum'
$ws.Range('T132').Value = 'Oxidative phosphorylation'
$ws.Range('J133').Value = 'MurNAc, chitobiose, fructose'
$ws.Range('L133').Value = 'methanol, formate'
$ws.Range('N133').Value = 'nitrilase'
$ws.Range('P133').Value = 'partial sulfur oxidation'
$ws.Range('Q133').Value = 'LPS export, phosphate, lipoprotein-release, phospholipid/cholesterol'
$ws.Range('T133').Value = 'Oxidative phosphorylation'
$ws.Range('J134').Value = 'MurNAc, chitobiose, fructose, glycolate, galacturonate'
$ws.Range('N134').Value = 'ammonia_assimilation'
$ws.Range('P134').Value = 'thiosulfate'
$ws.Range('Q134').Value = 'glutamate, amino acid/amide, glutamate/aspartate, heme, phosphate, phospholipid'
$ws.Range('T134').Value = 'Oxidative phosphorylation'
$ws.Range('J135').Value = 'carotenoid synthesis, carbon fixation (RuBisCo), carbon fixation (reductive TCA), '
$ws.Range('K135').Value = 'chitobiose, glucose, starch/glycogen, galacturonate'
$ws.Range('N135').Value = 'ammonia_assimilation, nitrogen fixation'
$ws.Range('P135').Value = 'sulfate_red_ass, alkane_sulfonate'
$ws.Range('Q135').Value = 'amino acid/amide, branched amino, cobalt/nickel, iron, LPS export, lipoprotein release, macrolide, molybdate, manganese/zinc/iron, phosphate, sulfate, sulfonate'
$ws.Range('S135').Value = 'one chemotaxis protein'
$ws.Range('T135').Value = 'Oxidative phosphorylation'
$ws.Range('K136').Value = 'can make dextrin?'
$ws.Range('N136').Value = 'ammonia_assimilation, nitronate monooxygenase'
$ws.Range('S136').Value = 'one type IV gene'
$ws.Range('K137').Value = 'Wood-Ljungdahl, formate'
$ws.Range('L137').Value = 'chitobiose, glucose, glyceraldehyde, starch/glycogen, galacturonate'
$ws.Range('N137').Value = 'ammonia_assimilation, partial denitrification, partial nitrification, nitrate/nitrite transporter'
$ws.Range('P137').Value = 'sulfate_red_ass, thiosulfate'
$ws.Range('Q137').Value = 'branched amino (lots), glutamate, amino acid/amide (lots), carbohydrate, dipeptide, glutathione, heme, iron, iron(III), LPS export, lipoprotein release, molybdate, monosaccharide, oligopeptide, phospholipid/cholesterol, phosphonate, putrescine, ribose, spermidine/putrescine, sulfate, tungstate, urea'
$ws.Range('T137').Value = 'Oxidative phosphorylation'
$ws.Range('J138').Value = 'carbon fixation (reductive TCA), formate, glucose, NAG, glycolate, starch/glycogen'
$ws.Range('N138').Value = 'ammonia_assimlation, nitrate_red_ass, 3/4 denitrification, nitrogen fixation, partial nitrate_red_ass'
$ws.Range('P138').Value = 'sulfate_red_ass, sulfate oxidation (SOX), thiosulfate'
$ws.Range('Q138').Value = 'LPS export, molybdate, phosphate, phospholipid/cholesterol, zinc'
$ws.Range('S138').Value = 'chemotaxis, flagellum'
$ws.Range('T138').Value = 'Oxidative phosphorylation'
$ws.Range('J139').Value = 'one carotenoid genes'
$ws.Range('K139').Value = 'fructose, glucose'
$ws.Range('N139').Value = 'ammonia_assimilation'
$ws.Range('P139').Value = 'thiosulfate'
$ws.Range('Q139').Value = 'ribose, xylose, amino acid/amide, biotin, branched amino, carbohydrate, iron, monosaccharide, sorbitol/mannitol'
$ws.Range('T139').Value = 'Oxidative phosphorylation'
$ws.Range('J140').Value = 'carbon fixation (RuBisCo), formate, MurNAc, fructose, glucose, sorbose, galactose, sucrose, glycolate, cellobiose, glucoside, starch/glycogen, trehalose, maltose, xylose'
$ws.Range('N140').Value = 'nitrate_red_ass, , nitrogen fixation, nitroalkane, ammonia_assimilation'
$ws.Range('P140').Value = 'sulfate_red_ass, sulfate oxidation (SOX), thiosulfate, alkane_sulfonate, methanesulfonate'
$ws.Range('Q140').Value = 'nitrate/sulfonate/bicarbonate, sugar, xylose, arginine, amino acid, amino acid/amide (lots), arginine/ornithine, biotin, branched amino, carbohydrate, dipeptide, extacellular solute, glucose/mannose, heme, iron, LPS export, lipoprotein release, maltose, microcin C, molybdate, molybdenum, monosaccharide, nitrate/nitrite, phopholipid/cholesterol, phosphonate, putrescine, raffinose/stachyose/melibiose, ribose, sorbitol, sulfate, tungstate, urea, type II secretion '
$ws.Range('S140').Value = 'chemotaxis (purine?), flagellum'
$ws.Range('T140').Value = 'Oxidative phosphorylation'

# Restore the sheet view: scrolled so column F is at the left edge of row
# 132, with K140 as the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 132
$win.ScrollColumn = 6
$ws.Range("K140").Select() | Out-Null

